$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on specific rows to reflect repulled/recalculated data
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F15").Value = -7
$ws.Range("F19").Value = 2
$ws.Range("F24").Value = -2
$ws.Range("F30").Value = -2
